$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update shared-string text: model holdings as-of date 2021-03-26 -> 2021-03-29
$ws.Cells.Item(58, 1).Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Update Weight (D) / Percent Change (E) columns for rows 2-55
$ws.Range("D2").Value = 0.01633780098557183
$ws.Range("E2").Value = -0.005755163511187544
$ws.Range("D3").Value = 0.04966638253609159
$ws.Range("E3").Value = 0.007765323407699087
$ws.Range("D4").Value = 0.0149733232960058
$ws.Range("E4").Value = -0.003238707994609569
$ws.Range("D5").Value = 0.009691033223098482
$ws.Range("E5").Value = 0.01356797420741529
$ws.Range("D6").Value = 0.01579089031160304
$ws.Range("E6").Value = 0.01478626489138057
$ws.Range("D7").Value = 0.02071956312238998
$ws.Range("E7").Value = 0.01565153193060165
$ws.Range("D8").Value = 0.004469019189716881
$ws.Range("E8").Value = -0.02824225122349089
$ws.Range("D9").Value = 0.006761721850812897
$ws.Range("E9").Value = -0.005371686016288235
$ws.Range("D10").Value = 0.01403224247371304
$ws.Range("E10").Value = -0.004341926729986456
$ws.Range("D11").Value = 0.008926191468857434
$ws.Range("E11").Value = 0.01298038357762699
$ws.Range("D12").Value = 0.01501302997431108
$ws.Range("E12").Value = -0.03277835587929212
$ws.Range("D13").Value = 0.002968562400184083
$ws.Range("E13").Value = -0.01650038372985418
$ws.Range("D14").Value = 0.006154209672742115
$ws.Range("E14").Value = -0.01599767306573596
$ws.Range("D15").Value = 0.01463812732890246
$ws.Range("E15").Value = -0.01553936424011848
$ws.Range("D16").Value = 0.01069444051174748
$ws.Range("E16").Value = -0.01972062448644196
$ws.Range("D17").Value = 0.0218429041073955
$ws.Range("E17").Value = -0.01267265656430194
$ws.Range("D18").Value = 0.008740286102890418
$ws.Range("E18").Value = 0.00203315608382848
$ws.Range("D19").Value = 0.0171772392206946
$ws.Range("E19").Value = 0.006730127933062491
$ws.Range("D20").Value = 0.012062368125822
$ws.Range("E20").Value = 0.005558253736981378
$ws.Range("D21").Value = 0.007299975167052831
$ws.Range("E21").Value = 0.008248102936324608
$ws.Range("D22").Value = 0.01384448195966127
$ws.Range("E22").Value = 0.004476040021063943
$ws.Range("D23").Value = 0.0200013278954712
$ws.Range("E23").Value = -0.003646588462721256
$ws.Range("D24").Value = 0.00996419564196423
$ws.Range("E24").Value = -0.01269299990527606
$ws.Range("D25").Value = 0.02099441795705234
$ws.Range("E25").Value = 0.0001612253123739471
$ws.Range("D26").Value = 0.01148102329968355
$ws.Range("E26").Value = 0.0009524943445646805
$ws.Range("D27").Value = 0.02091607863354349
$ws.Range("E27").Value = -0.02386046591239765
$ws.Range("D28").Value = 0.05522938580603418
$ws.Range("E28").Value = 0.001485025987954858
$ws.Range("D29").Value = 0.02133391005656249
$ws.Range("E29").Value = -0.000230361667818535
$ws.Range("D30").Value = 0.03145810408422242
$ws.Range("E30").Value = -0.01601556033769214
$ws.Range("D31").Value = 0.0160149075794516
$ws.Range("E31").Value = -0.01293330325588382
$ws.Range("D32").Value = 0.0135920516343041
$ws.Range("E32").Value = -0.007006369426751591
$ws.Range("D33").Value = 0.02061909220769457
$ws.Range("E33").Value = -0.01591089896579145
$ws.Range("D34").Value = 0.03953867349100791
$ws.Range("E34").Value = 0.01040138685158021
$ws.Range("D35").Value = 0.01143161777700534
$ws.Range("E35").Value = -0.007174581482746967
$ws.Range("D36").Value = 0.009657966022149167
$ws.Range("E36").Value = 0.007076806944706471
$ws.Range("D37").Value = 0.01163373127887074
$ws.Range("E37").Value = -0.007133864876206375
$ws.Range("D38").Value = 0.007418639551561889
$ws.Range("E38").Value = 0.005527770465912019
$ws.Range("D39").Value = 0.01170774192351517
$ws.Range("E39").Value = -0.009565667011375178
$ws.Range("D40").Value = 0.01813625314113999
$ws.Range("E40").Value = -0.0009511128019783488
$ws.Range("D41").Value = 0.01715269919164364
$ws.Range("E41").Value = -0.001009445525993136
$ws.Range("D42").Value = 0.03304506935812853
$ws.Range("E42").Value = -0.008480085096324363
$ws.Range("D43").Value = 0.01134484892097102
$ws.Range("E43").Value = -0.0043758857516627
$ws.Range("D44").Value = 0.02154389980117369
$ws.Range("E44").Value = 0.004589519037137446
$ws.Range("D45").Value = 0.01411982499118805
$ws.Range("E45").Value = -0.01954194672641285
$ws.Range("D46").Value = 0.00824095834997271
$ws.Range("E46").Value = -0.007535366460510229
$ws.Range("D47").Value = 0.0134449416458288
$ws.Range("E47").Value = 0.003137254901960818
$ws.Range("D48").Value = 0.009901804082692738
$ws.Range("E48").Value = -0.02218672346467876
$ws.Range("D49").Value = 0.01478064826685231
$ws.Range("E49").Value = 0.001268331351565566
$ws.Range("D50").Value = 0.008517342868143561
$ws.Range("E50").Value = 0.004352345069507768
$ws.Range("D51").Value = 0.01111233702768174
$ws.Range("E51").Value = -0.02203672787979949
$ws.Range("D52").Value = 0.008775566462958389
$ws.Range("E52").Value = 0.008240861618798778
$ws.Range("D53").Value = 0.1413883212244238
$ws.Range("D54").Value = 0.04369882679784191
$ws.Range("E54").Value = -0.001139528994682282
$ws.Range("E55").Value = -0.002652611804404903

$ws.Protect()
